$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the AutoFilter over the full data range (A1:D25), filtering
# column D ("Abweichung") to only show rows with value "Ja".
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
[void]$ws.Range("A1:D25").AutoFilter(4, @("Ja"), 7)

# Keep the _FilterDatabase defined name in sync with the new filter range
# (Excel normally does this automatically when AutoFilter is applied
# through the UI).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Tabelle1!_FilterDatabase") {
        $n.RefersTo = "=Tabelle1!`$A`$1:`$D`$25"
    }
}

# Move the active selection to A18 (reflects where the user ended up after
# filtering/importing the data).
[void]$ws.Range("A18").Select()

# Set up the page for printing (A4, portrait) which Excel records in a
# pageSetup element once print options are touched.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
